# Generate Report for Handoff
# Refresh the localization-status report: flip the in-flight rows from
# "In Translation" to "Ready for handoff" and stamp the new handoff / HO
# xliff generation timestamps, then re-fit the columns that now hold the
# longer "Ready for handoff" status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-09-01 19:08:55"

# --- zh-cn sheet ------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-09-01 19:08:51"

# --- de-de sheet ------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-09-01 19:08:55"

# Re-fit the status columns now that they hold the longer "Ready for
# handoff" text (was previously sized for "In Translation").
$newStatusColWidth = 16.25

$overview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$zhcn.Columns.Item(3).ColumnWidth = $newStatusColWidth
$dede.Columns.Item(3).ColumnWidth = $newStatusColWidth
